$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2021.0,4.0,25.0,7.0,7.0),
    @(2021.0,4.0,28.0,3.0,3.0),
    @(2021.0,4.0,29.0,4.0,4.0),
    @(2021.0,4.0,30.0,5.0,5.0),
    @(2021.0,4.0,31.0,-1.0,-1.0),
    @(2021.0,6.0,22.0,2.0,2.0),
    @(2021.0,6.0,28.0,1.0,1.0),
    @(2021.0,6.0,29.0,2.0,2.0),
    @(2021.0,6.0,30.0,3.0,3.0),
    @(2021.0,6.0,31.0,-1.0,-1.0),
    @(2021.0,2.0,15.0,1.0,1.0),
    @(2021.0,2.0,28.0,7.0,7.0),
    @(2021.0,2.0,29.0,-1.0,-1.0),
    @(2021.0,2.0,30.0,-1.0,-1.0),
    @(2021.0,2.0,31.0,-1.0,-1.0),
    @(2021.0,12.0,15.0,3.0,3.0),
    @(2021.0,12.0,28.0,2.0,2.0),
    @(2021.0,12.0,29.0,3.0,3.0),
    @(2021.0,12.0,30.0,4.0,4.0),
    @(2021.0,12.0,31.0,5.0,5.0),
    @(2022.0,4.0,25.0,1.0,1.0),
    @(2022.0,4.0,28.0,4.0,4.0),
    @(2022.0,4.0,29.0,5.0,5.0),
    @(2022.0,4.0,30.0,6.0,6.0),
    @(2022.0,4.0,31.0,-1.0,-1.0),
    @(2022.0,6.0,22.0,3.0,3.0),
    @(2022.0,6.0,28.0,2.0,2.0),
    @(2022.0,6.0,29.0,3.0,3.0),
    @(2022.0,6.0,30.0,4.0,4.0),
    @(2022.0,6.0,31.0,-1.0,-1.0),
    @(2022.0,2.0,15.0,2.0,2.0),
    @(2022.0,2.0,28.0,1.0,1.0),
    @(2022.0,2.0,29.0,-1.0,-1.0),
    @(2022.0,2.0,30.0,-1.0,-1.0),
    @(2022.0,2.0,31.0,-1.0,-1.0),
    @(2022.0,12.0,15.0,4.0,4.0),
    @(2022.0,12.0,28.0,3.0,3.0),
    @(2022.0,12.0,29.0,4.0,4.0),
    @(2022.0,12.0,30.0,5.0,5.0),
    @(2022.0,12.0,31.0,6.0,6.0)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
